$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Market Number" column (U) with header + data values
$ws.Range("U1").Value = "Market Number"
$ws.Range("U2").Value = 0
$ws.Range("U3").Value = 0

# Match the author's saved view state: active cell sits just below the
# new column's last data row, with the window scrolled so column D is
# the first visible column.
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$ws.Range("U4").Select()
